# "Generate Report for Handoff" - refresh the handoff-status report with the
# latest handoff timestamps for the f8c50ff7-66db-4c9b-846d-42ed8bcd99ac file
# (row 7 in each sheet).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest Handoff Date" column (D)
$wsOverview.Range("D7").Value = "2016-25-14 04:25:18"

# zh-cn sheet: "Latest Handoff Datetime" column (E)
$wsZhCn.Range("E7").Value = "2016-03-14 04:25:16"

# de-de sheet: "Latest Handoff Datetime" column (E)
$wsDeDe.Range("E7").Value = "2016-03-14 04:25:18"
